# Update database and shift quarterly periods:
# drop the oldest quarter (column D) and append the newest quarter's data
# (new last column, M) with updated figures from the latest release.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the oldest quarter column (D). This shifts E:M left to D:L,
#    automatically carrying over values, styles and shared strings.
$ws.Columns("D:D").Delete()

# 2. The newest quarter becomes the new column M. Give it the same
#    formatting as its neighbor (column L) by copying formats across.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)   # xlPasteFormats
$ws.Columns("M").ColumnWidth = 30.17      # matches OOXML width 31 (Q4 column)

# 3. Populate column M with the newest quarter's figures.
$ws.Cells.Item(11, 13).Value = 5358
$ws.Cells.Item(12, 13).Value = -4606
$ws.Cells.Item(13, 13).Value = 753
$ws.Cells.Item(14, 13).Value = -547
$ws.Cells.Item(16, 13).Value = 330
$ws.Cells.Item(17, 13).Value = 535
$ws.Cells.Item(18, 13).Value = -16
$ws.Cells.Item(19, 13).Value = 239
$ws.Cells.Item(20, 13).Value = 758
$ws.Cells.Item(21, 13).Value = 473
$ws.Cells.Item(22, 13).Value = 1231
$ws.Cells.Item(24, 13).Value = 1231
$ws.Cells.Item(26, 13).Value = 2420
